{"js": "const body = context.document.body;\n\n// 1) Update experience years in the Professional Summary: \"21 years\" -> \"15+ years\"\nconst searchResults = body.search(\"21 years of experience\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"15+ years of experience\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove the EDUCATION section entirely: the \"EDUCATION\" heading paragraph\n//    plus the two degree paragraphs that follow it.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet eduIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"EDUCATION\" && paragraphs.items[i].style === \"Heading 2\") {\n    eduIndex = i;\n    break;\n  }\n}\n\nif (eduIndex !== -1) {\n  // Delete in reverse order so indices stay valid as items are removed.\n  paragraphs.items[eduIndex + 2].delete();\n  paragraphs.items[eduIndex + 1].delete();\n  paragraphs.items[eduIndex].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update experience years in the Professional Summary: \"21 years\" -> \"15+ years\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"21 years of experience\"\n$find.Replacement.Text = \"15+ years of experience\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Remove the EDUCATION section entirely: the \"EDUCATION\" Heading2 paragraph\n#    plus the two Heading3 degree paragraphs that follow it.\n$eduHeading = $null\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text.TrimEnd(\"`r\")\n    if ($paraText -eq \"EDUCATION\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $eduHeading = $p\n        break\n    }\n}\n\nif ($eduHeading -ne $null) {\n    $degree1 = $eduHeading.Next()\n    $degree2 = $degree1.Next()\n    $rangeStart = $eduHeading.Range.Start\n    $rangeEnd = $degree2.Range.End\n    $killRange = $d.Range($rangeStart, $rangeEnd)\n    $killRange.Delete()\n}\n"}
